# Update benchmark data for row 9 (Gemini 2L), matching the committed
# diff for ros2_common_benchmark_data.xlsx.
#
# Columns (1-indexed):
#   H = color_delay_cur   (8)
#   I = color_delay_avg   (9)
#   J = color_delay_min   (10)
#   K = color_delay_max   (11)
#   P = depth_delay_cur   (16)
#   Q = depth_delay_avg   (17)
#   R = depth_delay_min   (18)
#   S = depth_delay_max   (19)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

# The source data stores these numeric-looking values as text, so force
# each cell's format to Text before assigning to preserve that type.
function Set-TextValue($ws, $row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextValue $ws $row 8  "108.58"   # H9: 73.12 -> 108.58
Set-TextValue $ws $row 9  "98.72"    # I9: 55.12 -> 98.72
Set-TextValue $ws $row 10 "30.06"    # J9: 0 -> 30.06

# K9: 378.99 -> "" (emptied out)
$ws.Cells.Item($row, 11).ClearContents()

Set-TextValue $ws $row 16 "114.4"    # P9: 66.2 -> 114.4
Set-TextValue $ws $row 17 "79.88"    # Q9: 42.5 -> 79.88
Set-TextValue $ws $row 18 "29.73"    # R9: 0.07 -> 29.73
Set-TextValue $ws $row 19 "155.01"   # S9: 148.47 -> 155.01
